$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AgricultureMapping")

# Insert a new row at row 11, shifting existing rows 11-29 down to 12-30
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 2).Value = "06_crude_oil_and_ngl"
$ws.Cells.Item(11, 3).Value = "06_02_natural_gas_liquids"
$ws.Cells.Item(11, 1).Value = "Natural Gas Liquids"

$ws.Range("D16").Select()
